$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.6058719545628528, 2.273965417999197, 13.7752703967488,  3.711505138989949, 3.70610558205782,  42),
    @(0.2253023586752266, 2.175657567196922, 13.37304018303112, 3.656916759106107, 3.695312707466266, 41),
    @(0.5126286334978025, 2.228379559222172, 14.40588458795522, 3.795508475547805, 3.808640173520903, 40),
    @(0.2827631271565333, 2.17940400322231,  13.49293633865711, 3.673273245847239, 3.710249930404507, 39),
    @(0.4944692049487463, 2.291119334397112, 15.22527009109669, 3.901957212873649, 3.922455350943079, 38),
    @(0.3870010927830395, 2.123771186146944, 13.74064264440705, 3.706837283238509, 3.737431878041005, 37),
    @(0.4119344969644807, 2.107164434821323, 14.21970027773409, 3.770901785744902, 3.801504848368828, 36),
    @(0.2440861792744094, 1.822874919107301, 12.15750528843192, 3.486761432681037, 3.528986910503495, 35),
    @(0.2707931634234895, 1.979019481310122, 13.02697252254846, 3.609289753199161, 3.653242154482546, 34),
    @(0.4558814406081306, 1.700108481721239, 12.37740496525853, 3.51815363013876,  3.542580466416441, 33)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
